{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst oldText1 = \"\u25cf Estabelecer em qual meio de comunica\u00e7\u00e3o ter\u00e3o as reuni\u00f5es virtuais com o cliente \u2013 Whatsapp, e-mail ou skype, pelo menos a cada uma etapa iremos realizar entrevistas presenciais;\";\nconst newText1 = \"\u25cf Quando for necess\u00e1rio entrar em contato com o cliente, definir antes quem ser\u00e1 o membro da equipe a realizar o contato para que n\u00e3o tenha repeti\u00e7\u00f5es de assuntos.\";\n\nconst oldText2 = \"\u25cf Duas vezes por semanas realizar check list do trabalho feito durante a semana e o que ser\u00e1 realizado posteriormente;\";\nconst newText2 = \"\u25cf Em caso de d\u00favidas, sempre encaminhar e-mail para o cliente, deixando em c\u00f3pia os outros membros da equipe.\";\n\nconst oldText3 = \"\u25cf Documentar todas as mudan\u00e7as e enviar para o cliente checar para termos feedbacks cont\u00ednuos do processo;\";\nconst oldText4 = \"\u25cf Caso algum Stakeholder n\u00e3o der conta de seu trabalho, faremos o trabalho em pares e revisar.\";\n\n// Collect the paragraphs that need to be removed entirely (delete after\n// updating the two that simply get new wording, so indices stay valid).\nconst toDelete = [];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  if (text === oldText1) {\n    para.insertText(newText1, \"Replace\");\n  } else if (text === oldText2) {\n    para.insertText(newText2, \"Replace\");\n  } else if (text === oldText3 || text === oldText4) {\n    toDelete.push(para);\n  } else if (text === \"\" && i > 0) {\n    // The bold blank paragraph that immediately follows the \"Caso algum\n    // Stakeholder\" bullet is also removed by the diff. Only the blank\n    // bold paragraph encountered right after that bullet is deleted; the\n    // remaining blank paragraphs elsewhere in the document stay intact.\n    const prevText = paragraphs.items[i - 1].text.trim();\n    if (prevText === oldText4) {\n      toDelete.push(para);\n    }\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Reword the first bullet about defining the communication channel.\n$old1 = \"\u25cf Estabelecer em qual meio de comunica\u00e7\u00e3o ter\u00e3o as reuni\u00f5es virtuais com o cliente \u2013 Whatsapp, e-mail ou skype, pelo menos a cada uma etapa iremos realizar entrevistas presenciais; \"\n$new1 = \"\u25cf Quando for necess\u00e1rio entrar em contato com o cliente, definir antes quem ser\u00e1 o membro da equipe a realizar o contato para que n\u00e3o tenha repeti\u00e7\u00f5es de assuntos.\"\n$find1 = $d.Content.Find\n$find1.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null\n\n# 2) Reword the second bullet about the biweekly checklist.\n$old2 = \"\u25cf Duas vezes por semanas realizar check list do trabalho feito durante a semana e o que ser\u00e1 realizado posteriormente; \"\n$new2 = \"\u25cf Em caso de d\u00favidas, sempre encaminhar e-mail para o cliente, deixando em c\u00f3pia os outros membros da equipe.\"\n$find2 = $d.Content.Find\n$find2.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null\n\n# 3) Remove the remaining two bullets (\"Documentar todas...\" and \"Caso algum\n#    Stakeholder...\") plus the blank bold paragraph that used to sit right\n#    after them. Walk backwards so deleting a paragraph never invalidates\n#    the index of paragraphs still to be inspected.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t.StartsWith(\"\u25cf Documentar todas as mudan\u00e7as\")) {\n    $p.Range.Delete()\n  }\n  elseif ($t.StartsWith(\"\u25cf Caso algum Stakeholder\")) {\n    # Delete this bullet together with the blank bold paragraph right after it.\n    $next = $d.Paragraphs.Item($i + 1)\n    $rng = $d.Range($p.Range.Start, $next.Range.End)\n    $rng.Delete()\n  }\n}\n"}
